$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks from A2 and A3
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()

# Update row 2 links to the new single URL (shared string now reused by both A2/B2)
$ws.Range("A2").Value = "https://www.proximity.mu/"
$ws.Range("B2").Value = "https://www.proximity.mu/"

# Clear row 3 contents entirely (keep cell styling)
$ws.Range("A3:B3").ClearContents()

# Move the active selection to B2
$ws.Range("B2").Select()
